$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2023.0588
$ws.Range("J17").Value = 2023.0588
$ws.Range("L17").Value = 6069.1764
$ws.Range("N17").Value = -6405.1764

$ws.Range("H42").Value = 131.33333
$ws.Range("I42").Value = 139.5
$ws.Range("J42").Value = 115
$ws.Range("K42").Value = 418.5
$ws.Range("L42").Value = 345
$ws.Range("M42").Value = -188.5
$ws.Range("N42").Value = -805

$ws.Range("H88").Value = 3057.5
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 3057.5
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 3057.5
$ws.Range("N88").Value = -3869.5
$ws.Range("M88").ClearContents()

$ws.Range("H91").Value = 3057.5
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 3057.5
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 3057.5
$ws.Range("N91").Value = -5865.5
$ws.Range("M91").ClearContents()

$ws.Range("H92").Value = 15625426
$ws.Range("I92").Value = 16667084
$ws.Range("K92").Value = 16667084
$ws.Range("M92").Value = -16665836

$ws.Range("H138").Value = 2585.5715
$ws.Range("I138").Value = 1314.3
$ws.Range("J138").Value = 2911.5386
$ws.Range("K138").Value = 3942.9
$ws.Range("L138").Value = 8734.6158
$ws.Range("M138").Value = 1197.1
$ws.Range("N138").Value = -19014.6158

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2217.6155
$ws.Range("J2").Value = 2645.1
$ws.Range("L2").Value = 2645.1
$ws.Range("N2").Value = -2871.1

$ws.Range("H32").Value = 17097362
$ws.Range("I32").Value = 17335340
$ws.Range("K32").Value = 17335340
$ws.Range("M32").Value = -17335053

$ws.Range("H97").Value = 1316.3889
$ws.Range("I97").Value = 1269.5714
$ws.Range("J97").Value = 1480.25
$ws.Range("K97").Value = 1269.5714
$ws.Range("L97").Value = 1480.25
$ws.Range("M97").Value = -773.5714
$ws.Range("N97").Value = -2472.25

$ws.Range("H102").Value = 1411.7059
$ws.Range("I102").Value = 1133.3334
$ws.Range("K102").Value = 1133.3334
$ws.Range("M102").Value = 488.6666

$ws.Range("H116").Value = 2217.6155
$ws.Range("J116").Value = 2645.1
$ws.Range("L116").Value = 2645.1
$ws.Range("N116").Value = -7233.1

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2217.6155
$ws.Range("J3").Value = 2645.1
$ws.Range("L3").Value = 2645.1
$ws.Range("N3").Value = -2873.1

$ws.Range("H86").Value = 3604.257
$ws.Range("I86").Value = 3454.48
$ws.Range("J86").Value = 3978.7
$ws.Range("K86").Value = 3454.48
$ws.Range("L86").Value = 3978.7
$ws.Range("M86").Value = -2331.48
$ws.Range("N86").Value = -6224.7

$ws.Range("H89").Value = 3604.257
$ws.Range("I89").Value = 3454.48
$ws.Range("J89").Value = 3978.7
$ws.Range("K89").Value = 17272.4
$ws.Range("L89").Value = 19893.5
$ws.Range("M89").Value = -11656.4
$ws.Range("N89").Value = -31125.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H28").Value = 48046.4
$ws.Range("J28").Value = 48046.4
$ws.Range("L28").Value = 48046.4
$ws.Range("N28").Value = -48536.4

$ws.Range("H68").Value = 49972.7
$ws.Range("J68").Value = 49972.7
$ws.Range("L68").Value = 49972.7
$ws.Range("N68").Value = -51470.7

$ws.Range("H71").Value = 49972.7
$ws.Range("J71").Value = 49972.7
$ws.Range("L71").Value = 149918.1
$ws.Range("N71").Value = -157406.1

$ws.Range("H99").Value = 3431.8333
$ws.Range("I99").Value = 3223.75
$ws.Range("K99").Value = 3223.75
$ws.Range("M99").Value = -1725.75

$ws.Range("H107").Value = 981.64703
$ws.Range("I107").Value = 400.1111
$ws.Range("J107").Value = 1635.875
$ws.Range("K107").Value = 400.1111
$ws.Range("L107").Value = 1635.875
$ws.Range("M107").Value = 1519.8889
$ws.Range("N107").Value = -5475.875

$ws.Range("H126").Value = 3431.8333
$ws.Range("I126").Value = 3223.75
$ws.Range("K126").Value = 9671.25
$ws.Range("M126").Value = -7201.25

$ws.Range("H132").Value = 1911.9048
$ws.Range("I132").Value = 1691.3055
$ws.Range("J132").Value = 3235.5
$ws.Range("K132").Value = 5073.916499999999
$ws.Range("L132").Value = 9706.5
$ws.Range("M132").Value = -2543.916499999999
$ws.Range("N132").Value = -14766.5

$ws.Range("H133").Value = 29887
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H20").Value = 2000
$ws.Range("I20").Value = 2000
$ws.Range("K20").Value = 6000
$ws.Range("M20").Value = -5773

$ws.Range("H35").Value = 3735.3333
$ws.Range("I35").Value = 3735.3333
$ws.Range("K35").Value = 11205.9999
$ws.Range("M35").Value = -10917.9999

$ws.Range("H40").Value = 526.8570999999999
$ws.Range("I40").Value = 362.8
$ws.Range("J40").Value = 937
$ws.Range("K40").Value = 1451.2
$ws.Range("L40").Value = 3748
$ws.Range("M40").Value = -1382.2
$ws.Range("N40").Value = -3886

$ws.Range("H41").Value = 700
$ws.Range("J41").Value = 700
$ws.Range("L41").Value = 2100
$ws.Range("N41").Value = -2776

$ws.Range("H56").Value = 7369.4287
$ws.Range("I56").Value = 7369.4287
$ws.Range("K56").Value = 7369.4287
$ws.Range("M56").Value = -6839.4287

$ws.Range("H64").Value = 394
$ws.Range("I64").Value = 394
$ws.Range("K64").Value = 1182
$ws.Range("M64").Value = -912

$ws.Range("H67").Value = 394
$ws.Range("I67").Value = 394
$ws.Range("K67").Value = 1182
$ws.Range("M67").Value = -246

$ws.Range("H81").Value = 3691.2222
$ws.Range("I81").Value = 1110.5
$ws.Range("K81").Value = 3331.5
$ws.Range("M81").Value = -2208.5

$ws.Range("H84").Value = 3691.2222
$ws.Range("I84").Value = 1110.5
$ws.Range("K84").Value = 9994.5
$ws.Range("M84").Value = -4378.5

$ws.Range("H86").Value = 400
$ws.Range("I86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("M86").ClearContents()

$ws.Range("H88").Value = 3814.8333
$ws.Range("J88").Value = 3814.8333
$ws.Range("L88").Value = 11444.4999
$ws.Range("N88").Value = -12300.4999

$ws.Range("H89").Value = 400
$ws.Range("I89").Value = 0
$ws.Range("K89").Value = 0
$ws.Range("M89").ClearContents()

$ws.Range("H91").Value = 3814.8333
$ws.Range("J91").Value = 3814.8333
$ws.Range("L91").Value = 11444.4999
$ws.Range("N91").Value = -14408.4999

$ws.Range("H117").Value = 224278
$ws.Range("J117").Value = 251937.75
$ws.Range("L117").Value = 755813.25
$ws.Range("N117").Value = -762697.25

$ws.Range("H128").Value = 1979899
$ws.Range("I128").Value = 1979899
$ws.Range("K128").Value = 5939697
$ws.Range("M128").Value = -5934717

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 12322.739
$ws.Range("I113").Value = 8592.904
$ws.Range("K113").Value = 8592.904
$ws.Range("M113").Value = -6422.904

$ws.Range("H121").Value = 80602
$ws.Range("J121").Value = 80602
$ws.Range("L121").Value = 80602
$ws.Range("N121").Value = -84096

$ws.Range("H126").Value = 2234.0667
$ws.Range("J126").Value = 2693.625
$ws.Range("L126").Value = 8080.875
$ws.Range("N126").Value = -13020.875

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 10913.392
$ws.Range("J2").Value = 10913.392
$ws.Range("L2").Value = 10913.392
$ws.Range("N2").Value = -11137.392

$ws.Range("H12").Value = 1597.4
$ws.Range("J12").Value = 1746.75
$ws.Range("L12").Value = 1746.75
$ws.Range("N12").Value = -2086.75

$ws.Range("H93").Value = 31251664
$ws.Range("I93").Value = 62501452
$ws.Range("J93").Value = 1873.3125
$ws.Range("K93").Value = 62501452
$ws.Range("L93").Value = 1873.3125
$ws.Range("M93").Value = -62500204
$ws.Range("N93").Value = -4369.3125

$ws.Range("H132").Value = 3498
$ws.Range("I132").Value = 3498
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 10494
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -7964
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value = 750000.5
$ws.Range("J5").Value = 750000.5
$ws.Range("L5").Value = 750000.5
$ws.Range("N5").Value = -750224.5

$ws.Range("H126").Value = 1735
$ws.Range("I126").Value = 1916.5714
$ws.Range("K126").Value = 5749.7142
$ws.Range("M126").Value = -3279.7142
